$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item("Table1")

# Reorder the grade rows: David/Eve (previously rows 5-6) move up to rows 3-4,
# Bob/Cecilia (previously rows 3-4) move down to rows 5-6.
$ws.Range("A3").Value = 45678
$ws.Range("B3").Value = "David"
$ws.Range("C3").Value = 4
$ws.Range("D3").Formula = "=ROUND(4.2,0)"

$ws.Range("A4").Value = 56789
$ws.Range("B4").Value = "Eve"
$ws.Range("C4").Value = 5
$ws.Range("D4").Formula = "=ROUND(4.9, 0)"

$ws.Range("A5").Value = 23456
$ws.Range("B5").Value = "Bob"
$ws.Range("C5").Value = 2
$ws.Range("D5").Formula = "=ROUND(2.2, 0)"

$ws.Range("A6").Value = 34567
$ws.Range("B6").Value = "Cecilia"
$ws.Range("C6").Value = 3
$ws.Range("D6").Formula = "=ROUND(3.2,0)"

# Grow the table by one row so Table1 / its AutoFilter cover A1:D8.
$lo.ListRows.Add()

# Move Frank's record down into the newly added last row (row 8).
$ws.Range("A8").Value = 67890
$ws.Range("B8").Value = "Frank"
$ws.Range("C8").Value = "H"
$ws.Range("D8").Value = "H"

# Row 7 used to hold Frank; remove those leftover cells entirely ...
$ws.Range("C7").Clear()
$ws.Range("D7").Clear()

# ... then populate it with the new student (only Student number + Name).
$ws.Range("A7").Value = 99999
$ws.Range("B7").Value = "Empty"

# Match the saved selection from the edit.
$ws.Range("B9").Select() | Out-Null
